$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.676.01'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.850.34'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.24%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '262.65'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.65%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5376'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.46%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3191'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.40%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07019'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.05'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7770'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07844'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.89%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.855.60'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.58%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '89.60'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.056'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.89%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.17'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.71%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008015'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.57%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.693.66'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.087.04'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.649'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.051'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.418'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.15'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.220'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.41%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.13'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.64%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.71'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.321'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08759'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.116'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04881'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7377'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.66%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.145'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.88%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.117'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.355'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +7.44%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01752'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4843'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9082'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.51'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.34%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.927'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.17%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.757'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.82%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4202'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.161'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.55%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.17'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05837'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8995'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.95%  '
